$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update countries & provincias Spain
# Columns: A=Pais, B=Casos totales, C=Nuevos casos, D=Casos activos,
#          E=Recuperados, F=Casos criticos, G=Muertes hoy, H=Muertes
# The data rows stay sorted by column B (Casos totales) descending; some
# countries' stats were refreshed which shuffled a few rows' positions.

# Polonia (row 33): refreshed Casos activos / Recuperados
$ws.Cells.Item(33, 1).Value = "Polonia"
$ws.Cells.Item(33, 2).Value = 10169
$ws.Cells.Item(33, 3).Value = 0
$ws.Cells.Item(33, 4).Value = 1740
$ws.Cells.Item(33, 5).Value = 8003
$ws.Cells.Item(33, 6).Value = 160
$ws.Cells.Item(33, 7).Value = 0
$ws.Cells.Item(33, 8).Value = 426

# Rumania (row 35): refreshed Recuperados / Muertes hoy / Muertes
$ws.Cells.Item(35, 1).Value = "Rumania"
$ws.Cells.Item(35, 2).Value = 9710
$ws.Cells.Item(35, 3).Value = 0
$ws.Cells.Item(35, 4).Value = 2406
$ws.Cells.Item(35, 5).Value = 6777
$ws.Cells.Item(35, 6).Value = 288
$ws.Cells.Item(35, 7).Value = 3
$ws.Cells.Item(35, 8).Value = 527

# Oman moves up to row 68 with fresh stats; Uzbekistan/Irak shift down
$ws.Cells.Item(68, 1).Value = "Oman"
$ws.Cells.Item(68, 2).Value = 1716
$ws.Cells.Item(68, 3).Value = 102
$ws.Cells.Item(68, 4).Value = 307
$ws.Cells.Item(68, 5).Value = 1401
$ws.Cells.Item(68, 6).Value = 3
$ws.Cells.Item(68, 7).Value = 0
$ws.Cells.Item(68, 8).Value = 8

$ws.Cells.Item(69, 1).Value = "Uzbekistan"
$ws.Cells.Item(69, 2).Value = 1716
$ws.Cells.Item(69, 3).Value = 0
$ws.Cells.Item(69, 4).Value = 450
$ws.Cells.Item(69, 5).Value = 1259
$ws.Cells.Item(69, 6).Value = 8
$ws.Cells.Item(69, 7).Value = 0
$ws.Cells.Item(69, 8).Value = 7

$ws.Cells.Item(70, 1).Value = "Irak"
$ws.Cells.Item(70, 2).Value = 1631
$ws.Cells.Item(70, 3).Value = 0
$ws.Cells.Item(70, 4).Value = 1146
$ws.Cells.Item(70, 5).Value = 402
$ws.Cells.Item(70, 6).Value = 0
$ws.Cells.Item(70, 7).Value = 0
$ws.Cells.Item(70, 8).Value = 83

# Armenia moves up to row 72 with fresh stats; Azerbaiyan shifts down
$ws.Cells.Item(72, 1).Value = "Armenia"
$ws.Cells.Item(72, 2).Value = 1523
$ws.Cells.Item(72, 3).Value = 50
$ws.Cells.Item(72, 4).Value = 659
$ws.Cells.Item(72, 5).Value = 840
$ws.Cells.Item(72, 6).Value = 30
$ws.Cells.Item(72, 7).Value = 0
$ws.Cells.Item(72, 8).Value = 24

$ws.Cells.Item(73, 1).Value = "Azerbaiyan"
$ws.Cells.Item(73, 2).Value = 1518
$ws.Cells.Item(73, 3).Value = 0
$ws.Cells.Item(73, 4).Value = 907
$ws.Cells.Item(73, 5).Value = 591
$ws.Cells.Item(73, 6).Value = 14
$ws.Cells.Item(73, 7).Value = 0
$ws.Cells.Item(73, 8).Value = 20

# Islas Malvinas moves up to row 195 with fresh stats; Montserrat shifts down
$ws.Cells.Item(195, 1).Value = "Islas Malvinas"
$ws.Cells.Item(195, 2).Value = 12
$ws.Cells.Item(195, 3).Value = 1
$ws.Cells.Item(195, 4).Value = 11
$ws.Cells.Item(195, 5).Value = 1
$ws.Cells.Item(195, 6).Value = 0
$ws.Cells.Item(195, 7).Value = 0
$ws.Cells.Item(195, 8).Value = 0

$ws.Cells.Item(196, 1).Value = "Montserrat"
$ws.Cells.Item(196, 2).Value = 11
$ws.Cells.Item(196, 3).Value = 0
$ws.Cells.Item(196, 4).Value = 2
$ws.Cells.Item(196, 5).Value = 9
$ws.Cells.Item(196, 6).Value = 1
$ws.Cells.Item(196, 7).Value = 0
$ws.Cells.Item(196, 8).Value = 0

# Update the "last refreshed" timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 23 de Abril de 2020 a las 09:22"
